# Rever_DailyTrack_BALRAJ_2022.xlsx — "Add files via upload"
#
# Fills in the two blank daily-entries (rows 31 & 32) on the FEB-22 sheet
# with the next log entry (No. 21, dated 25-Feb-2022) and moves the saved
# cell selection to D38.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")

# --- Row 31: new entry #21, 25-Feb-2022, RPA GSS, completed -----------------
$ws.Cells.Item(31, 1).Value = 21                 # A31 - No

# B31 - Date: copy the number format from an existing date cell (B29) first,
# so the new cell keeps/reuses the workbook's date style instead of minting
# a new one, then write the serial date value (25-Feb-2022 = 44617).
$ws.Range("B29").Copy()
$ws.Range("B31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(31, 2).Value = 44617

$ws.Cells.Item(31, 3).Value = "RPA GSS"           # C31 - Application
$ws.Cells.Item(31, 4).Value = "1. Mr Kabilan san has completed the captcha issue at token system based on the sony captcha code( whereas I have supported to that logical approach and testing), it is tested and running smoothly"  # D31 - Task

# E31 - % of completion: reuse the percentage format from E29 the same way.
$ws.Range("E29").Copy()
$ws.Range("E31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(31, 5).Value = 1

$ws.Cells.Item(31, 6).Value = "Completed"         # F31 - Status

$ws.Rows.Item(31).RowHeight = 28.8                # two-line wrapped task text

# --- Row 32: continuation line, invoice generation task, in progress -------
$ws.Cells.Item(32, 4).Value = "2. Task of invoice generation has been completed, tested and it is updating at master file success, whereas the integration and  formatting the data is work in progress"  # D32 - Task

# E32 - % of completion (80%): reuse the percentage format from E29.
$ws.Range("E29").Copy()
$ws.Range("E32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(32, 5).Value = 0.8

$ws.Cells.Item(32, 6).Value = "WIP"               # F32 - Status

# --- Update the saved selection on the sheet -------------------------------
$ws.Range("D38").Select() | Out-Null
